$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1687.9375
$ws.Range("I19").Value = 623.8333
$ws.Range("J19").Value = 2326.4
$ws.Range("K19").Value = 623.8333
$ws.Range("L19").Value = 2326.4
$ws.Range("M19").Value = -448.8333
$ws.Range("N19").Value = -2676.4

$ws.Range("H62").Value = 4079.6365
$ws.Range("I62").Value = 1556.4445
$ws.Range("J62").Value = 15434
$ws.Range("K62").Value = 1556.4445
$ws.Range("L62").Value = 15434
$ws.Range("M62").Value = -932.4445000000001
$ws.Range("N62").Value = -16682

$ws.Range("H65").Value = 4079.6365
$ws.Range("I65").Value = 1556.4445
$ws.Range("J65").Value = 15434
$ws.Range("K65").Value = 7782.2225
$ws.Range("L65").Value = 77170
$ws.Range("M65").Value = -4662.2225
$ws.Range("N65").Value = -83410

$ws.Range("H131").Value = 1642.7858
$ws.Range("I131").Value = 499.91666
$ws.Range("J131").Value = 8500
$ws.Range("K131").Value = 1499.74998
$ws.Range("L131").Value = 25500
$ws.Range("M131").Value = 3540.25002
$ws.Range("N131").Value = -35580

$ws.Range("H137").Value = 1853.0435
$ws.Range("I137").Value = 1212.5
$ws.Range("J137").Value = 2079.1177
$ws.Range("K137").Value = 3637.5
$ws.Range("L137").Value = 6237.353099999999
$ws.Range("M137").Value = -1087.5
$ws.Range("N137").Value = -11337.3531

$ws.Range("H138").Value = 4527.559
$ws.Range("I138").Value = 4242.1875
$ws.Range("J138").Value = 4633.744
$ws.Range("K138").Value = 12726.5625
$ws.Range("L138").Value = 13901.232
$ws.Range("M138").Value = -7586.5625
$ws.Range("N138").Value = -24181.232

$ws.Range("H139").Value = 70000
$ws.Range("J139").Value = 70000
$ws.Range("L139").Value = 70000
$ws.Range("N139").Value = -80280

$ws.Range("H140").Value = 163932.67
$ws.Range("J140").Value = 163932.67
$ws.Range("L140").Value = 163932.67
$ws.Range("N140").Value = -174292.67

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2845.19
$ws.Range("I32").Value = 2203.7
$ws.Range("K32").Value = 2203.7
$ws.Range("M32").Value = -1916.7

$ws.Range("H46").Value = 12055.2
$ws.Range("I46").Value = 10092
$ws.Range("J46").Value = 15000
$ws.Range("K46").Value = 10092
$ws.Range("L46").Value = 15000
$ws.Range("M46").Value = -9773
$ws.Range("N46").Value = -15638

$ws.Range("H61").Value = 2546.182
$ws.Range("I61").Value = 1089.1333
$ws.Range("J61").Value = 5668.4287
$ws.Range("K61").Value = 1089.1333
$ws.Range("L61").Value = 5668.4287
$ws.Range("M61").Value = -877.1333
$ws.Range("N61").Value = -6092.4287

$ws.Range("H74").Value = 1685.1428
$ws.Range("I74").Value = 1119.2
$ws.Range("K74").Value = 1119.2
$ws.Range("M74").Value = -245.2

$ws.Range("H77").Value = 1685.1428
$ws.Range("I77").Value = 1119.2
$ws.Range("K77").Value = 5596
$ws.Range("M77").Value = -1228

$ws.Range("H122").Value = 1589.875
$ws.Range("I122").Value = 1549.963
$ws.Range("K122").Value = 4649.889
$ws.Range("M122").Value = -2199.889

$ws.Range("H132").Value = 2330.4
$ws.Range("I132").Value = 1750.037
$ws.Range("J132").Value = 3200.9443
$ws.Range("K132").Value = 5250.111
$ws.Range("L132").Value = 9602.832900000001
$ws.Range("M132").Value = -2720.111
$ws.Range("N132").Value = -14662.8329

$ws.Range("H136").Value = 2546.182
$ws.Range("I136").Value = 1089.1333
$ws.Range("J136").Value = 5668.4287
$ws.Range("K136").Value = 3267.3999
$ws.Range("L136").Value = 17005.2861
$ws.Range("M136").Value = -717.3998999999999
$ws.Range("N136").Value = -22105.2861

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2513.0625
$ws.Range("I31").Value = 1169.2858
$ws.Range("J31").Value = 3558.2222
$ws.Range("K31").Value = 1169.2858
$ws.Range("L31").Value = 3558.2222
$ws.Range("M31").Value = -874.2858000000001
$ws.Range("N31").Value = -4148.2222

$ws.Range("H34").Value = 2513.0625
$ws.Range("I34").Value = 1169.2858
$ws.Range("J34").Value = 3558.2222
$ws.Range("K34").Value = 1169.2858
$ws.Range("L34").Value = 3558.2222
$ws.Range("M34").Value = -967.2858000000001
$ws.Range("N34").Value = -3962.2222

$ws.Range("H58").Value = 1451117
$ws.Range("I58").Value = 1976977.8
$ws.Range("J58").Value = 5000
$ws.Range("K58").Value = 1976977.8
$ws.Range("L58").Value = 5000
$ws.Range("M58").Value = -1976774.8
$ws.Range("N58").Value = -5406

$ws.Range("H134").Value = 2406.1667
$ws.Range("I134").Value = 1455.7916
$ws.Range("J134").Value = 6207.6665
$ws.Range("K134").Value = 4367.3748
$ws.Range("L134").Value = 18622.9995
$ws.Range("M134").Value = -1832.3748
$ws.Range("N134").Value = -23692.9995

$ws.Range("H136").Value = 1451117
$ws.Range("I136").Value = 1976977.8
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 5930933.4
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -5928383.4
$ws.Range("N136").Value = -20100

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").Value = $null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 2723.36
$ws.Range("I107").Value = 2443
$ws.Range("J107").Value = 2784.9023
$ws.Range("K107").Value = 7329
$ws.Range("L107").Value = 8354.706900000001
$ws.Range("M107").Value = -5409
$ws.Range("N107").Value = -12194.7069

$ws.Range("H113").Value = 1214.7391
$ws.Range("J113").Value = 811.381
$ws.Range("L113").Value = 2434.143
$ws.Range("N113").Value = -6774.143

$ws.Range("H132").Value = 2044.3846
$ws.Range("I132").Value = 861
$ws.Range("J132").Value = 2399.4
$ws.Range("K132").Value = 7749
$ws.Range("L132").Value = 21594.6
$ws.Range("M132").Value = -5219
$ws.Range("N132").Value = -26654.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 702.7
$ws.Range("J107").Value = 1204.2
$ws.Range("L107").Value = 1204.2
$ws.Range("N107").Value = -5044.2

$ws.Range("H113").Value = 1381.4615
$ws.Range("I113").Value = 1275.5714
$ws.Range("J113").Value = 1505
$ws.Range("K113").Value = 1275.5714
$ws.Range("L113").Value = 1505
$ws.Range("M113").Value = 894.4286
$ws.Range("N113").Value = -5845

$ws.Range("H122").Value = 1943.76
$ws.Range("I122").Value = 1852.75
$ws.Range("J122").Value = 2105.5557
$ws.Range("K122").Value = 5558.25
$ws.Range("L122").Value = 6316.6671
$ws.Range("M122").Value = -3108.25
$ws.Range("N122").Value = -11216.6671

$ws.Range("I132").Value = 4809601
$ws.Range("J132").Value = 5191
$ws.Range("K132").Value = 14428803
$ws.Range("L132").Value = 15573
$ws.Range("M132").Value = -14426273
$ws.Range("N132").Value = -20633

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2758.625
$ws.Range("I22").Value = 775
$ws.Range("K22").Value = 775
$ws.Range("M22").Value = -480

$ws.Range("H27").Value = 2758.625
$ws.Range("I27").Value = 775
$ws.Range("K27").Value = 775
$ws.Range("M27").Value = -668

$ws.Range("H68").Value = 2987.182
$ws.Range("I68").Value = 2607.375
$ws.Range("K68").Value = 2607.375
$ws.Range("M68").Value = -1858.375

$ws.Range("H71").Value = 2987.182
$ws.Range("I71").Value = 2607.375
$ws.Range("K71").Value = 13036.875
$ws.Range("M71").Value = -9292.875

$ws.Range("H122").Value = 5369.6665
$ws.Range("I122").Value = 3768.7856
$ws.Range("J122").Value = 8571.429
$ws.Range("K122").Value = 11306.3568
$ws.Range("L122").Value = 25714.287
$ws.Range("M122").Value = -8856.356800000001
$ws.Range("N122").Value = -30614.287

$ws.Range("H136").Value = 4422
$ws.Range("I136").Value = 3023.7896
$ws.Range("J136").Value = 6319.5713
$ws.Range("K136").Value = 9071.3688
$ws.Range("L136").Value = 18958.7139
$ws.Range("M136").Value = -6521.3688
$ws.Range("N136").Value = -24058.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 1450
$ws.Range("I58").Value = 1450
$ws.Range("K58").Value = 1450
$ws.Range("M58").Value = -1142

$ws.Range("H107").Value = 842.5
$ws.Range("I107").Value = 828.7143
$ws.Range("J107").Value = 856.2857
$ws.Range("K107").Value = 2486.1429
$ws.Range("L107").Value = 2568.8571
$ws.Range("M107").Value = -566.1428999999998
$ws.Range("N107").Value = -6408.8571

$ws.Range("H132").Value = 891.475
$ws.Range("I132").Value = 544.48
$ws.Range("J132").Value = 1469.8
$ws.Range("K132").Value = 1633.44
$ws.Range("L132").Value = 4409.4
$ws.Range("M132").Value = 896.5599999999999
$ws.Range("N132").Value = -9469.4

$ws.Range("H136").Value = 15875839
$ws.Range("I136").Value = 27780762
$ws.Range("J136").Value = 2609.2666
$ws.Range("K136").Value = 83342286
$ws.Range("L136").Value = 7827.7998
$ws.Range("M136").Value = -83339736
$ws.Range("N136").Value = -12927.7998
